$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor refresh of the existing last-row timestamp (new data pull recomputed it).
$ws.Range("A6").Value = 45806.40665445602

# Append the new price-check row, matching the date-column formatting
# already used by the rows above it.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("A7").Value = 45807.39290911525
$ws.Range("B7").Value = "EVOWHEY PROTEIN"
$ws.Range("C7").Value = "2Kg"
$ws.Range("D7").Value = "37,90€"
